$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the header row (row 1) labels:
# B1: kitchens_1 -> living_rooms_1
# C1: living_rooms_1 -> kitchens_1
# D1: bedrooms_2 -> kitchens_2
# E1: kitchens_2 -> living_rooms_2
# F1: living_rooms_2 -> bedrooms_2
$ws.Range("B1").Value = "living_rooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "living_rooms_2"
$ws.Range("F1").Value = "bedrooms_2"
